$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sorted (descending by value) language data, with "Swedish" and "Uzbek" removed
$ws.Cells.Item(2,1).Value = "English"
$ws.Cells.Item(2,2).Value = 21.72716039894293
$ws.Cells.Item(3,1).Value = "Chinese"
$ws.Cells.Item(3,2).Value = 16.80593093581323
$ws.Cells.Item(4,1).Value = "Spanish"
$ws.Cells.Item(4,2).Value = 6.918744996955311
$ws.Cells.Item(5,1).Value = "Arabic"
$ws.Cells.Item(5,2).Value = 5.136289736146665
$ws.Cells.Item(6,1).Value = "Japanese"
$ws.Cells.Item(6,2).Value = 4.750037304829528
$ws.Cells.Item(7,1).Value = "German"
$ws.Cells.Item(7,2).Value = 4.288633200707078
$ws.Cells.Item(8,1).Value = "Russian"
$ws.Cells.Item(8,2).Value = 3.70910415200979
$ws.Cells.Item(9,1).Value = "Portuguese"
$ws.Cells.Item(9,2).Value = 3.458853659834345
$ws.Cells.Item(10,1).Value = "Malay-Indonesian"
$ws.Cells.Item(10,2).Value = 3.061707035891857
$ws.Cells.Item(11,1).Value = "French"
$ws.Cells.Item(11,2).Value = 2.671806180422946
$ws.Cells.Item(12,1).Value = "Italian"
$ws.Cells.Item(12,2).Value = 2.069087834841883
$ws.Cells.Item(13,1).Value = "Korean"
$ws.Cells.Item(13,2).Value = 1.633516881765565
$ws.Cells.Item(14,1).Value = "Turkish"
$ws.Cells.Item(14,2).Value = 1.611540257882197
$ws.Cells.Item(15,1).Value = "Dutch"
$ws.Cells.Item(15,2).Value = 1.253039309003411
$ws.Cells.Item(16,1).Value = "Persian"
$ws.Cells.Item(16,2).Value = 1.243315222853494
$ws.Cells.Item(17,1).Value = "Thai"
$ws.Cells.Item(17,2).Value = 0.9931697281063218
$ws.Cells.Item(18,1).Value = "Polish"
$ws.Cells.Item(18,2).Value = 0.8751043766664499
$ws.Cells.Item(19,1).Value = "Urdu"
$ws.Cells.Item(19,2).Value = 0.8356045333349099
$ws.Cells.Item(20,1).Value = "Bengali"
$ws.Cells.Item(20,2).Value = 0.5770005971137531
$ws.Cells.Item(21,1).Value = "Vietnamese"
$ws.Cells.Item(21,2).Value = 0.5741921523881593

# Remove the two trailing rows that no longer exist (Uzbek, Vietnamese old position)
$ws.Range("A22:B23").Delete()
